$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 165572860
$ws.Range("I11").Value = 165572860
$ws.Range("K11").Value = 165572860
$ws.Range("M11").Value = -165572720
$ws.Range("H17").Value = 18199.166
$ws.Range("J17").Value = 21479
$ws.Range("L17").Value = 64437
$ws.Range("N17").Value = -64773
$ws.Range("H28").Value = 2438
$ws.Range("I28").Value = 2837.2
$ws.Range("J28").Value = 442
$ws.Range("K28").Value = 2837.2
$ws.Range("L28").Value = 442
$ws.Range("M28").Value = -2352.2
$ws.Range("N28").Value = -1412
$ws.Range("H32").Value = 2350
$ws.Range("J32").Value = 3000
$ws.Range("L32").Value = 3000
$ws.Range("N32").Value = -3652
$ws.Range("H33").Value = 418.20587
$ws.Range("I33").Value = 238.11539
$ws.Range("J33").Value = 1003.5
$ws.Range("K33").Value = 238.11539
$ws.Range("L33").Value = 1003.5
$ws.Range("M33").Value = -9.115389999999991
$ws.Range("N33").Value = -1461.5
$ws.Range("H53").Value = 1389.3182
$ws.Range("I53").Value = 887.5
$ws.Range("K53").Value = 887.5
$ws.Range("M53").Value = -250.5
$ws.Range("H58").Value = 2794.375
$ws.Range("I58").Value = 2143.3333
$ws.Range("J58").Value = 4747.5
$ws.Range("K58").Value = 6429.999899999999
$ws.Range("L58").Value = 14242.5
$ws.Range("M58").Value = -6279.999899999999
$ws.Range("N58").Value = -14542.5
$ws.Range("H62").Value = 8334.375
$ws.Range("J62").Value = 8668.75
$ws.Range("L62").Value = 8668.75
$ws.Range("N62").Value = -9916.75
$ws.Range("H64").Value = 5759.4287
$ws.Range("I64").Value = 8001.5
$ws.Range("J64").Value = 5523.421
$ws.Range("K64").Value = 8001.5
$ws.Range("L64").Value = 5523.421
$ws.Range("M64").Value = -7753.5
$ws.Range("N64").Value = -6019.421
$ws.Range("H65").Value = 8334.375
$ws.Range("J65").Value = 8668.75
$ws.Range("L65").Value = 43343.75
$ws.Range("N65").Value = -49583.75
$ws.Range("H67").Value = 5759.4287
$ws.Range("I67").Value = 8001.5
$ws.Range("J67").Value = 5523.421
$ws.Range("K67").Value = 8001.5
$ws.Range("L67").Value = 5523.421
$ws.Range("M67").Value = -7143.5
$ws.Range("N67").Value = -7239.421
$ws.Range("H81").Value = 132000
$ws.Range("J81").Value = 132000
$ws.Range("L81").Value = 132000
$ws.Range("N81").Value = -133996
$ws.Range("H84").Value = 132000
$ws.Range("J84").Value = 132000
$ws.Range("L84").Value = 396000
$ws.Range("N84").Value = -405984
$ws.Range("H98").Value = 1441.275
$ws.Range("I98").Value = 1468.4324
$ws.Range("K98").Value = 1468.4324
$ws.Range("M98").Value = 29.56760000000008
$ws.Range("H112").Value = 2969.139
$ws.Range("J112").Value = 3626.6155
$ws.Range("L112").Value = 10879.8465
$ws.Range("N112").Value = -13095.8465
$ws.Range("H113").Value = 2231.4285
$ws.Range("I113").Value = 1967.1
$ws.Range("J113").Value = 2892.25
$ws.Range("K113").Value = 1967.1
$ws.Range("L113").Value = 2892.25
$ws.Range("M113").Value = 1286.9
$ws.Range("N113").Value = -9400.25
$ws.Range("H122").Value = 1441.275
$ws.Range("I122").Value = 1468.4324
$ws.Range("K122").Value = 4405.2972
$ws.Range("M122").Value = -1955.2972
$ws.Range("H132").Value = 3361.5173
$ws.Range("I132").Value = 3515.6738
$ws.Range("K132").Value = 10547.0214
$ws.Range("M132").Value = -8017.0214
$ws.Range("H137").Value = 43439.9
$ws.Range("I137").Value = 53681.312
$ws.Range("J137").Value = 2474.25
$ws.Range("K137").Value = 161043.936
$ws.Range("L137").Value = 7422.75
$ws.Range("M137").Value = -158493.936
$ws.Range("N137").Value = -12522.75
$ws.Range("H138").Value = 3548.8235
$ws.Range("I138").Value = 2037.3243
$ws.Range("J138").Value = 7543.5
$ws.Range("K138").Value = 6111.9729
$ws.Range("L138").Value = 22630.5
$ws.Range("M138").Value = -971.9728999999998
$ws.Range("N138").Value = -32910.5
$ws.Range("H141").Value = 8357.571
$ws.Range("I141").Value = 8434
$ws.Range("J141").Value = 7899
$ws.Range("K141").Value = 25302
$ws.Range("L141").Value = 23697
$ws.Range("M141").Value = -20122
$ws.Range("N141").Value = -34057

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1735.3334
$ws.Range("I2").Value = 1315.7333
$ws.Range("K2").Value = 1315.7333
$ws.Range("M2").Value = -1202.7333
$ws.Range("H32").Value = 6616464.5
$ws.Range("I32").Value = 3145638.8
$ws.Range("K32").Value = 3145638.8
$ws.Range("M32").Value = -3145351.8
$ws.Range("H45").Value = 3192.7144
$ws.Range("I45").Value = 2890.75
$ws.Range("J45").Value = 4159
$ws.Range("K45").Value = 2890.75
$ws.Range("L45").Value = 4159
$ws.Range("M45").Value = -2513.75
$ws.Range("N45").Value = -4913
$ws.Range("H61").Value = 4812.375
$ws.Range("I61").Value = 4500
$ws.Range("J61").Value = 4999.8
$ws.Range("K61").Value = 4500
$ws.Range("L61").Value = 4999.8
$ws.Range("M61").Value = -4288
$ws.Range("N61").Value = -5423.8
$ws.Range("H74").Value = 1832.64
$ws.Range("I74").Value = 1765.3043
$ws.Range("K74").Value = 1765.3043
$ws.Range("M74").Value = -891.3043
$ws.Range("H77").Value = 1832.64
$ws.Range("I77").Value = 1765.3043
$ws.Range("K77").Value = 8826.521500000001
$ws.Range("M77").Value = -4458.521500000001
$ws.Range("H92").Value = 68229.25
$ws.Range("J92").Value = 68229.25
$ws.Range("L92").Value = 68229.25
$ws.Range("N92").Value = -73221.25
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H110").Value = 3003.6667
$ws.Range("I110").Value = 1211
$ws.Range("K110").Value = 1211
$ws.Range("M110").Value = 834
$ws.Range("H116").Value = 1735.3334
$ws.Range("I116").Value = 1315.7333
$ws.Range("K116").Value = 1315.7333
$ws.Range("M116").Value = 978.2666999999999
$ws.Range("H122").Value = 4663.6665
$ws.Range("I122").Value = 3912.7097
$ws.Range("K122").Value = 11738.1291
$ws.Range("M122").Value = -9288.1291
$ws.Range("H132").Value = 2564.2273
$ws.Range("I132").Value = 1770.2449
$ws.Range("J132").Value = 3561.795
$ws.Range("K132").Value = 5310.7347
$ws.Range("L132").Value = 10685.385
$ws.Range("M132").Value = -2780.7347
$ws.Range("N132").Value = -15745.385
$ws.Range("H134").Value = 21000
$ws.Range("J134").Value = 21000
$ws.Range("L134").Value = 21000
$ws.Range("N134").Value = -31140
$ws.Range("H136").Value = 4812.375
$ws.Range("I136").Value = 4500
$ws.Range("J136").Value = 4999.8
$ws.Range("K136").Value = 13500
$ws.Range("L136").Value = 14999.4
$ws.Range("M136").Value = -10950
$ws.Range("N136").Value = -20099.4
$ws.Range("H137").Value = 149777.44
$ws.Range("I137").Value = 112499.25
$ws.Range("K137").Value = 112499.25
$ws.Range("M137").Value = -107399.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1735.3334
$ws.Range("I3").Value = 1315.7333
$ws.Range("K3").Value = 1315.7333
$ws.Range("M3").Value = -1201.7333
$ws.Range("H13").Value = 78702.664
$ws.Range("J13").Value = 78702.664
$ws.Range("L13").Value = 78702.664
$ws.Range("N13").Value = -79038.664
$ws.Range("H20").Value = 16055.6
$ws.Range("I20").Value = 21617.809
$ws.Range("K20").Value = 21617.809
$ws.Range("M20").Value = -21370.809
$ws.Range("H81").Value = 346069.34
$ws.Range("I81").Value = 18709
$ws.Range("J81").Value = 509749.5
$ws.Range("K81").Value = 18709
$ws.Range("L81").Value = 509749.5
$ws.Range("M81").Value = -17648
$ws.Range("N81").Value = -511871.5
$ws.Range("H84").Value = 346069.34
$ws.Range("I84").Value = 18709
$ws.Range("J84").Value = 509749.5
$ws.Range("K84").Value = 56127
$ws.Range("L84").Value = 1529248.5
$ws.Range("M84").Value = -50823
$ws.Range("N84").Value = -1539856.5
$ws.Range("H86").Value = 5678.9546
$ws.Range("I86").Value = 4839.9375
$ws.Range("K86").Value = 4839.9375
$ws.Range("M86").Value = -3716.9375
$ws.Range("H89").Value = 5678.9546
$ws.Range("I89").Value = 4839.9375
$ws.Range("K89").Value = 24199.6875
$ws.Range("M89").Value = -18583.6875
$ws.Range("H94").Value = 1447.9474
$ws.Range("J94").Value = 1097.25
$ws.Range("L94").Value = 1097.25
$ws.Range("N94").Value = -1999.25
$ws.Range("H105").Value = 2930.7778
$ws.Range("I105").Value = 2229.5
$ws.Range("J105").Value = 4333.3335
$ws.Range("K105").Value = 2229.5
$ws.Range("L105").Value = 4333.3335
$ws.Range("M105").Value = -482.5
$ws.Range("N105").Value = -7827.3335
$ws.Range("H107").Value = 6499.25
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 2000
$ws.Range("M107").Value = -80
$ws.Range("H134").Value = 12651903
$ws.Range("I134").Value = 3248246.2
$ws.Range("J134").Value = 33339948
$ws.Range("K134").Value = 9744738.600000001
$ws.Range("L134").Value = 100019844
$ws.Range("M134").Value = -9742203.600000001
$ws.Range("N134").Value = -100024914
$ws.Range("H135").Value = 63249.75
$ws.Range("J135").Value = 63249.75
$ws.Range("L135").Value = 63249.75
$ws.Range("N135").Value = -73389.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3178.6
$ws.Range("I16").Value = 3123.25
$ws.Range("K16").Value = 3123.25
$ws.Range("M16").Value = -2836.25
$ws.Range("H31").Value = 4893.089
$ws.Range("I31").Value = 2421.9312
$ws.Range("J31").Value = 9372.0625
$ws.Range("K31").Value = 2421.9312
$ws.Range("L31").Value = 9372.0625
$ws.Range("M31").Value = -2126.9312
$ws.Range("N31").Value = -9962.0625
$ws.Range("H34").Value = 4893.089
$ws.Range("I34").Value = 2421.9312
$ws.Range("J34").Value = 9372.0625
$ws.Range("K34").Value = 2421.9312
$ws.Range("L34").Value = 9372.0625
$ws.Range("M34").Value = -2219.9312
$ws.Range("N34").Value = -9776.0625
$ws.Range("H74").Value = 46312.5
$ws.Range("I74").Value = 46312.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 46312.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -45438.5
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 46312.5
$ws.Range("I77").Value = 46312.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 138937.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -134569.5
$ws.Range("N77").ClearContents()
$ws.Range("H97").Value = 96354
$ws.Range("J97").Value = 96354
$ws.Range("L97").Value = 96354
$ws.Range("N97").Value = -98336
$ws.Range("H99").Value = 8355
$ws.Range("I99").Value = 1174.5
$ws.Range("J99").Value = 14099.4
$ws.Range("K99").Value = 1174.5
$ws.Range("L99").Value = 14099.4
$ws.Range("M99").Value = 323.5
$ws.Range("N99").Value = -17095.4
$ws.Range("H105").Value = 3960
$ws.Range("I105").Value = 3950
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 3950
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -2203
$ws.Range("N105").Value = -7494
$ws.Range("H113").Value = 3178.6
$ws.Range("I113").Value = 3123.25
$ws.Range("K113").Value = 3123.25
$ws.Range("M113").Value = -953.25
$ws.Range("H122").Value = 2152.0386
$ws.Range("I122").Value = 1917.7
$ws.Range("K122").Value = 5753.1
$ws.Range("M122").Value = -3303.1
$ws.Range("H126").Value = 8355
$ws.Range("I126").Value = 1174.5
$ws.Range("J126").Value = 14099.4
$ws.Range("K126").Value = 3523.5
$ws.Range("L126").Value = 42298.2
$ws.Range("M126").Value = -1053.5
$ws.Range("N126").Value = -47238.2
$ws.Range("H132").Value = 2273.652
$ws.Range("I132").Value = 2074.6
$ws.Range("K132").Value = 6223.799999999999
$ws.Range("M132").Value = -3693.799999999999
$ws.Range("H134").Value = 2029.8462
$ws.Range("I134").Value = 1505.1212
$ws.Range("K134").Value = 4515.363600000001
$ws.Range("M134").Value = -1980.363600000001
$ws.Range("H138").Value = 94998.336
$ws.Range("J138").Value = 94998.336
$ws.Range("L138").Value = 94998.336
$ws.Range("N138").Value = -105278.336
$ws.Range("H141").Value = 540420.4399999999
$ws.Range("J141").Value = 574622.2
$ws.Range("L141").Value = 574622.2
$ws.Range("N141").Value = -584982.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 51331996
$ws.Range("I4").Value = 40775630
$ws.Range("J4").Value = 75083816
$ws.Range("K4").Value = 122326890
$ws.Range("L4").Value = 225251448
$ws.Range("M4").Value = -122326778
$ws.Range("N4").Value = -225251672
$ws.Range("H24").Value = 799.6
$ws.Range("I24").Value = 799
$ws.Range("K24").Value = 2397
$ws.Range("M24").Value = -2167
$ws.Range("H69").Value = 398
$ws.Range("I69").Value = 398
$ws.Range("K69").Value = 1194
$ws.Range("M69").Value = -383
$ws.Range("H72").Value = 398
$ws.Range("I72").Value = 398
$ws.Range("K72").Value = 3582
$ws.Range("M72").Value = 474
$ws.Range("H107").Value = 889.5
$ws.Range("J107").Value = 905.0625
$ws.Range("L107").Value = 2715.1875
$ws.Range("N107").Value = -6555.1875
$ws.Range("H122").Value = 2121.923
$ws.Range("J122").Value = 2073.625
$ws.Range("L122").Value = 18662.625
$ws.Range("N122").Value = -23562.625
$ws.Range("H131").Value = 2052.6
$ws.Range("J131").Value = 2314.5557
$ws.Range("L131").Value = 6943.6671
$ws.Range("N131").Value = -17023.6671
$ws.Range("H133").Value = 2500
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 2500
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 7500
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -17620
$ws.Range("H137").Value = 2440.75
$ws.Range("J137").Value = 2188.6667
$ws.Range("L137").Value = 6566.000100000001
$ws.Range("N137").Value = -16766.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 16799.1
$ws.Range("J20").Value = 21996.8
$ws.Range("L20").Value = 21996.8
$ws.Range("N20").Value = -22486.8
$ws.Range("H80").Value = 2532.6667
$ws.Range("J80").Value = 2532.6667
$ws.Range("L80").Value = 2532.6667
$ws.Range("N80").Value = -4528.6667
$ws.Range("H83").Value = 2532.6667
$ws.Range("J83").Value = 2532.6667
$ws.Range("L83").Value = 12663.3335
$ws.Range("N83").Value = -22647.3335
$ws.Range("H97").Value = 766.0952
$ws.Range("J97").Value = 4000
$ws.Range("L97").Value = 4000
$ws.Range("N97").Value = -4992
$ws.Range("H102").Value = 2439.4783
$ws.Range("I102").Value = 2179.3684
$ws.Range("K102").Value = 2179.3684
$ws.Range("M102").Value = -557.3683999999998
$ws.Range("H107").Value = 1903.2858
$ws.Range("I107").Value = 1705.75
$ws.Range("J107").Value = 2166.6667
$ws.Range("K107").Value = 1705.75
$ws.Range("L107").Value = 2166.6667
$ws.Range("M107").Value = 214.25
$ws.Range("N107").Value = -6006.6667
$ws.Range("H122").Value = 4681.2856
$ws.Range("I122").Value = 4052.4
$ws.Range("J122").Value = 6253.5
$ws.Range("K122").Value = 12157.2
$ws.Range("L122").Value = 18760.5
$ws.Range("M122").Value = -9707.200000000001
$ws.Range("N122").Value = -23660.5
$ws.Range("H126").Value = 2149.4167
$ws.Range("I126").Value = 1974.125
$ws.Range("K126").Value = 5922.375
$ws.Range("M126").Value = -3452.375
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 2064.4167
$ws.Range("I132").Value = 1583.4117
$ws.Range("J132").Value = 3232.5715
$ws.Range("K132").Value = 4750.2351
$ws.Range("L132").Value = 9697.7145
$ws.Range("M132").Value = -2220.2351
$ws.Range("N132").Value = -14757.7145
$ws.Range("H136").Value = 30352.158
$ws.Range("J136").Value = 30352.158
$ws.Range("L136").Value = 91056.474
$ws.Range("N136").Value = -96156.474
$ws.Range("H138").Value = 79999.5
$ws.Range("J138").Value = 79999.5
$ws.Range("L138").Value = 79999.5
$ws.Range("N138").Value = -90279.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3030
$ws.Range("I7").Value = 2325
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 2325
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -2213
$ws.Range("N7").Value = -3724
$ws.Range("H46").Value = 3293.6333
$ws.Range("I46").Value = 1329.25
$ws.Range("J46").Value = 4603.222
$ws.Range("K46").Value = 1329.25
$ws.Range("L46").Value = 4603.222
$ws.Range("M46").Value = -1141.25
$ws.Range("N46").Value = -4979.222
$ws.Range("H58").Value = 11750
$ws.Range("J58").Value = 11750
$ws.Range("L58").Value = 11750
$ws.Range("N58").Value = -12270
$ws.Range("H61").Value = 9166.666999999999
$ws.Range("I61").Value = 11250
$ws.Range("K61").Value = 11250
$ws.Range("M61").Value = -11048
$ws.Range("H68").Value = 6000
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 6000
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H99").Value = 69824
$ws.Range("I99").Value = 42499.75
$ws.Range("J99").Value = 106256.336
$ws.Range("K99").Value = 42499.75
$ws.Range("L99").Value = 106256.336
$ws.Range("M99").Value = -39504.75
$ws.Range("N99").Value = -112246.336
$ws.Range("H100").Value = 2896
$ws.Range("I100").Value = 2346.25
$ws.Range("K100").Value = 2346.25
$ws.Range("M100").Value = -1805.25
$ws.Range("H113").Value = 9166.666999999999
$ws.Range("I113").Value = 11250
$ws.Range("K113").Value = 11250
$ws.Range("M113").Value = -9080
$ws.Range("H126").Value = 3030
$ws.Range("I126").Value = 2325
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 6975
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -4505
$ws.Range("N126").Value = -15440
$ws.Range("H132").Value = 4097.3
$ws.Range("I132").Value = 3899.4
$ws.Range("K132").Value = 11698.2
$ws.Range("M132").Value = -9168.200000000001
$ws.Range("H134").Value = 49499.5
$ws.Range("J134").Value = 49499.5
$ws.Range("L134").Value = 49499.5
$ws.Range("N134").Value = -59639.5
$ws.Range("H136").Value = 6627.3
$ws.Range("I136").Value = 5816.3335
$ws.Range("J136").Value = 6974.857
$ws.Range("K136").Value = 17449.0005
$ws.Range("L136").Value = 20924.571
$ws.Range("M136").Value = -14899.0005
$ws.Range("N136").Value = -26024.571
$ws.Range("H139").Value = 94233.336
$ws.Range("J139").Value = 94233.336
$ws.Range("L139").Value = 94233.336
$ws.Range("N139").Value = -104513.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6749.25
$ws.Range("J62").Value = 6749.25
$ws.Range("L62").Value = 6749.25
$ws.Range("N62").Value = -7997.25
$ws.Range("H65").Value = 6749.25
$ws.Range("J65").Value = 6749.25
$ws.Range("L65").Value = 33746.25
$ws.Range("N65").Value = -39986.25
$ws.Range("H81").Value = 4222.0713
$ws.Range("I81").Value = 2050.5
$ws.Range("K81").Value = 4101
$ws.Range("M81").Value = -3040
$ws.Range("H84").Value = 4222.0713
$ws.Range("I84").Value = 2050.5
$ws.Range("K84").Value = 20505
$ws.Range("M84").Value = -15201
$ws.Range("H113").Value = 352.26666
$ws.Range("I113").Value = 268.875
$ws.Range("J113").Value = 447.57144
$ws.Range("K113").Value = 806.625
$ws.Range("L113").Value = 1342.71432
$ws.Range("M113").Value = 1363.375
$ws.Range("N113").Value = -5682.71432
$ws.Range("H122").Value = 1630.3334
$ws.Range("I122").Value = 1656.96
$ws.Range("J122").Value = 1297.5
$ws.Range("K122").Value = 4970.88
$ws.Range("L122").Value = 3892.5
$ws.Range("M122").Value = -2520.88
$ws.Range("N122").Value = -8792.5
$ws.Range("H126").Value = 7632.7856
$ws.Range("I126").Value = 7425.7144
$ws.Range("J126").Value = 7839.857
$ws.Range("K126").Value = 22277.1432
$ws.Range("L126").Value = 23519.571
$ws.Range("M126").Value = -19807.1432
$ws.Range("N126").Value = -28459.571
$ws.Range("H132").Value = 4369.9507
$ws.Range("I132").Value = 3941.739
$ws.Range("J132").Value = 5683.1333
$ws.Range("K132").Value = 11825.217
$ws.Range("L132").Value = 17049.3999
$ws.Range("M132").Value = -9295.217000000001
$ws.Range("N132").Value = -22109.3999
$ws.Range("H136").Value = 32826.938
$ws.Range("I136").Value = 1877
$ws.Range("J136").Value = 94726.82000000001
$ws.Range("K136").Value = 5631
$ws.Range("L136").Value = 284180.46
$ws.Range("M136").Value = -3081
$ws.Range("N136").Value = -289280.46
